$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$text)
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $origStyle
}

$ws.Range('D2').Value = '64.901.11'
$ws.Range('E2').Value = '  -0.06%  '
$ws.Range('D3').Value = '3.517.30'
$ws.Range('E3').Value = '  -0.20%  '
Set-TextValue $ws.Range('D4') '0.999'
$ws.Range('E4').Value = '  -0.02%  '
Set-TextValue $ws.Range('D5') '596.08'
$ws.Range('E5').Value = '  +0.37%  '
Set-TextValue $ws.Range('D6') '134.56'
$ws.Range('E6').Value = '  -1.20%  '
$ws.Range('D7').Value = '3.517.83'
$ws.Range('E7').Value = '  -0.24%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('E9').Value = '  +0.59%  '
$ws.Range('E10').Value = '  +1.06%  '
Set-TextValue $ws.Range('D11') '7.16'
$ws.Range('E11').Value = '  +4.57%  '
Set-TextValue $ws.Range('D12') '0.383'
$ws.Range('E12').Value = '  -0.37%  '
$ws.Range('D13').Value = '4.115.70'
$ws.Range('E13').Value = '  +0.07%  '
Set-TextValue $ws.Range('D14') '27.25'
$ws.Range('E14').Value = '  +0.33%  '
$ws.Range('E15').Value = '  +0.22%  '
$ws.Range('B16').Value = 'TRON'
$ws.Range('C16').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue $ws.Range('D16') '0.117'
$ws.Range('E16').Value = '  -0.16%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.512.56'
$ws.Range('E17').Value = '  -0.33%  '
$ws.Range('D18').Value = '63.920.95'
$ws.Range('E18').Value = '  -1.48%  '
Set-TextValue $ws.Range('D19') '9.85'
$ws.Range('E19').Value = '  -1.84%  '
Set-TextValue $ws.Range('D20') '14.42'
$ws.Range('E20').Value = '  +1.82%  '
Set-TextValue $ws.Range('D21') '5.71'
$ws.Range('E21').Value = '  -1.84%  '
Set-TextValue $ws.Range('D22') '388.83'
$ws.Range('E22').Value = '  +0.11%  '
$ws.Range('E23').Value = '  +0.95%  '
$ws.Range('D24').Value = '3.656.08'
$ws.Range('E24').Value = '  -0.14%  '
Set-TextValue $ws.Range('D25') '74.22'
$ws.Range('E25').Value = '  +0.69%  '
$ws.Range('E26').Value = '  +0.25%  '
$ws.Range('E27').Value = '  +1.56%  '
$ws.Range('B28').Value = 'Fetch.AI'
$ws.Range('C28').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws.Range('D28') '1.62'
$ws.Range('E28').Value = '  +19.15%  '
$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D29') '7.81'
$ws.Range('E29').Value = '  +2.38%  '
$ws.Range('E30').Value = '  -0.16%  '
Set-TextValue $ws.Range('D31') '2.29'
$ws.Range('E31').Value = '  +1.76%  '
Set-TextValue $ws.Range('D32') '8.42'
$ws.Range('E32').Value = '  +3.10%  '
$ws.Range('D33').Value = '3.519.44'
$ws.Range('E33').Value = '  -0.40%  '
$ws.Range('E34').Value = '  +1.57%  '
$ws.Range('E35').Value = '  +0.01%  '
Set-TextValue $ws.Range('D36') '0.145'
$ws.Range('E36').Value = '  +1.52%  '
Set-TextValue $ws.Range('D37') '5.27'
$ws.Range('E37').Value = '  +6.70%  '
Set-TextValue $ws.Range('D38') '1.59'
$ws.Range('E38').Value = '  +1.83%  '
Set-TextValue $ws.Range('D39') '169.83'
$ws.Range('E39').Value = '  +0.51%  '
Set-TextValue $ws.Range('D40') '6.86'
$ws.Range('E40').Value = '  +0.28%  '
$ws.Range('E41').Value = '  +3.42%  '
Set-TextValue $ws.Range('D42') '0.823'
$ws.Range('E42').Value = '  +0.62%  '
$ws.Range('B43').Value = 'OKB'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range('D43') '42.58'
$ws.Range('E43').Value = '  +0.64%  '
$ws.Range('B44').Value = 'ONDO'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
Set-TextValue $ws.Range('D44') '1.24'
$ws.Range('E44').Value = '  +3.75%  '
$ws.Range('E45').Value = '  +0.28%  '
Set-TextValue $ws.Range('D46') '25.32'
$ws.Range('E46').Value = '  -3.30%  '
$ws.Range('E47').Value = '  +0.73%  '
$ws.Range('E48').Value = '  -0.43%  '
$ws.Range('E49').Value = '  +1.48%  '
$ws.Range('D50').Value = '2.372.79'
$ws.Range('E50').Value = '  -0.96%  '
Set-TextValue $ws.Range('D51') '0.897'
$ws.Range('E51').Value = '  +7.09%  '
